$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A100").Value = "Golang Architect / Principal Backend Architect Only Local to GA"
$ws.Range("B100").Value = "https://www.dice.com/job-detail/ac26ccd2-8c43-4730-bd1b-6494576424db"
$ws.Range("C100").Value = "Atlanta, Georgia"
$ws.Range("D100").Value = "Third Party"
$ws.Range("E100").Value = "Depends on Experience"
$ws.Range("F100").Value = "Dahl Consulting"

$ws.Range("A101").Value = "Remote, Lead - Integration/GoLang Developer (.NET/Python/GoLang)"
$ws.Range("B101").Value = "https://www.dice.com/job-detail/8d571df4-7128-4b87-b352-e24305235af8"
$ws.Range("C101").Value = "Remote"
$ws.Range("D101").Value = "Contract, Third Party"
$ws.Range("E101").Value = "Depends on Experience"
$ws.Range("F101").Value = "Swanktek"
